$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '25.792.17'
Set-TextValue 'E2' '  -1.30%  '
Set-TextValue 'D3' '1.632.33'
Set-TextValue 'E3' '  -1.32%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  -0.50%  '
Set-TextValue 'D5' '214.31'
Set-TextValue 'E5' '  -0.52%  '
Set-TextValue 'D6' '0.5015'
Set-TextValue 'E6' '  -1.59%  '
Set-TextValue 'D7' '1.001'
Set-TextValue 'E7' '  -0.51%  '
Set-TextValue 'D8' '0.2560'
Set-TextValue 'E8' '  -0.69%  '
Set-TextValue 'D9' '0.06360'
Set-TextValue 'E9' '  -0.69%  '
Set-TextValue 'D10' '19.63'
Set-TextValue 'E10' '  -1.39%  '
Set-TextValue 'D11' '0.07687'
Set-TextValue 'E11' '  -1.41%  '
Set-TextValue 'B12' 'WrappedEther'
Set-TextValue 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.651.60'
Set-TextValue 'E12' '  -0.25%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '4.247'
Set-TextValue 'E13' '  -0.74%  '
Set-TextValue 'B14' 'WrappedliquidstakedEther2.0'
Set-TextValue 'C14' 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue 'D14' '1.857.87'
Set-TextValue 'E14' '  -1.29%  '
Set-TextValue 'D15' '0.5412'
Set-TextValue 'E15' '  -1.75%  '
Set-TextValue 'D16' '0.0₅7901'
Set-TextValue 'E16' '  -1.29%  '
Set-TextValue 'D17' '63.49'
Set-TextValue 'E17' '  -0.68%  '
Set-TextValue 'D18' '25.812.77'
Set-TextValue 'E18' '  -1.31%  '
Set-TextValue 'D19' '1.001'
Set-TextValue 'D20' '200.48'
Set-TextValue 'E20' '  -4.15%  '
Set-TextValue 'D21' '4.320'
Set-TextValue 'E21' '  -1.89%  '
Set-TextValue 'D22' '9.884'
Set-TextValue 'E22' '  -1.64%  '
Set-TextValue 'D23' '5.932'
Set-TextValue 'E23' '  -1.59%  '
Set-TextValue 'E24' '  -0.38%  '
Set-TextValue 'E25' '  +10.96%  '
Set-TextValue 'D26' '141.37'
Set-TextValue 'E26' '  -1.60%  '
Set-TextValue 'D27' '0.1135'
Set-TextValue 'E27' '  -3.45%  '
Set-TextValue 'D28' '15.61'
Set-TextValue 'E28' '  -1.25%  '
Set-TextValue 'D29' '6.695'
Set-TextValue 'E29' '  -3.99%  '
Set-TextValue 'D30' '1.238'
Set-TextValue 'E30' '  -0.32%  '
Set-TextValue 'D31' '0.04979'
Set-TextValue 'E31' '  -2.42%  '
Set-TextValue 'D32' '3.263'
Set-TextValue 'E32' '  -2.44%  '
Set-TextValue 'D33' '3.182'
Set-TextValue 'E33' '  -1.14%  '
Set-TextValue 'D34' '1.536'
Set-TextValue 'E34' '  -1.93%  '
Set-TextValue 'E35' '  +0.11%  '
Set-TextValue 'D36' '1.168.09'
Set-TextValue 'E36' '  +0.36%  '
Set-TextValue 'D37' '2.625'
Set-TextValue 'E37' '  -4.73%  '
Set-TextValue 'D38' '0.8895'
Set-TextValue 'E38' '  -4.12%  '
Set-TextValue 'D39' '0.5563'
Set-TextValue 'E39' '  -2.17%  '
Set-TextValue 'D40' '0.01555'
Set-TextValue 'E40' '  -2.02%  '
Set-TextValue 'E41' '  -0.50%  '
Set-TextValue 'D42' '5.672'
Set-TextValue 'E42' '  +0.43%  '
Set-TextValue 'D43' '0.8035'
Set-TextValue 'E43' '  -3.57%  '
Set-TextValue 'D44' '99.33'
Set-TextValue 'E44' '  -1.05%  '
Set-TextValue 'D45' '1.770.28'
Set-TextValue 'E45' '  -1.25%  '
Set-TextValue 'E46' '  -0.84%  '
Set-TextValue 'D47' '0.4512'
Set-TextValue 'E47' '  -0.82%  '
Set-TextValue 'D48' '1.002'
Set-TextValue 'E48' '  -0.65%  '
Set-TextValue 'D49' '54.59'
Set-TextValue 'E49' '  -2.00%  '
Set-TextValue 'D50' '0.05071'
Set-TextValue 'E50' '  +0.56%  '
Set-TextValue 'E51' '  -0.40%  '
